# Journal de travail: log two new work entries and extend the table /
# totals row accordingly (mirrors dragging the table's resize handle
# down a few extra blank rows before filling in the latest entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")
$tbl = $ws.ListObjects.Item("Tableau1")

# --- New data row 1: 5/4/2023 - Implémentation - 2h ---------------------
$ws.Range("A30").NumberFormat = "m/d/yy"
$ws.Range("A30").Value = 45050
$ws.Range("B30").Value = "Implémentation"
$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = "Travail sur l'API, édition configuration de module"

# --- New data row 2: 5/4/2023 - Rédaction ------------------------------
$ws.Range("A31").NumberFormat = "m/d/yy"
$ws.Range("A31").Value = 45050
$ws.Range("B31").Value = "Rédaction"
$ws.Range("C31").NumberFormat = "General"
$ws.Range("D31").Value = "Rapport"

# --- Extra blank rows (32-38) carrying the table's column formatting ---
for ($r = 32; $r -le 37; $r++) {
    $ws.Range("A$r").NumberFormat = "m/d/yy"
    $ws.Range("C$r").NumberFormat = "General"
}
$ws.Range("C38").NumberFormat = "General"

# --- Totals row moves down to row 39 ------------------------------------
$ws.Range("C39").Formula = "=SUM(C2:C37)"

# Grow the table to cover the new rows (totals row included)
$tbl.Resize($ws.Range("A1:E39"))

# --- Restore view state (scroll position / selection) -------------------
$ws.Range("D31").Select()
